$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Work plan")

# Row 36: task got a start date (C36) and was marked done (D36), matching
# the "Good" styled date cells used elsewhere in the log (e.g. D43/C42).
$ws.Range("C42").Copy()
$ws.Range("C36").PasteSpecial(-4122)
$ws.Range("C36").Value = 43832

$ws.Range("D43").Copy()
$ws.Range("D36").PasteSpecial(-4122)
$ws.Range("D36").Value = 44167

# Row 44: task marked done (D44) with notes about the DFD rework.
$ws.Range("D43").Copy()
$ws.Range("D44").PasteSpecial(-4122)
$ws.Range("D44").Value = 44107

$ws.Range("F44").Value = "Final data flowdiagram created"
$ws.Range("G44").Value = "Update the specifications"
$ws.Range("H44").Value = "Had to rework the DFD from the ground"

# Row 49: task marked done (D49, text style matching D37) with notes about
# the production-ready push and writing up the report.
$ws.Range("D37").Copy()
$ws.Range("D49").PasteSpecial(-4122)
$ws.Range("D49").Value = "14/3/2020"

$ws.Range("F49").Value = "Production ready code Pushed to github"
$ws.Range("G49").Value = "Write report about it"

# Update view: scroll/select down near the bottom of the log.
$ws.Range("E52").Select()

$wb.Save()
